# Daily attendance processing - 2026-01-31 19:59:58
# Swap the order of recorder names in the "Recorded By" (column G) cells
# from "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$used.Replace(
    $oldText,
    $newText,
    1,      # xlWhole
    1,      # xlByRows (LookAt uses whole/part, this param is SearchOrder but kept for compat)
    $false, # MatchCase
    $false, # MatchByte
    $false  # SearchFormat
)
